$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily COVID overview rows for the United Kingdom (2021-11-08 .. 2021-11-14),
# appended after the existing last row (453).
$rows = @(
    @(454, "2021-11-08", 9333891, 32322,  57, 141862),
    @(455, "2021-11-09", 9366676, 33117, 262, 142124),
    @(456, "2021-11-10", 9406001, 39329, 214, 142338),
    @(457, "2021-11-11", 9448402, 42408, 195, 142533),
    @(458, "2021-11-12", 9487302, 40375, 145, 142678),
    @(459, "2021-11-13", 9524971, 38351, 157, 142835),
    @(460, "2021-11-14", 9561099, 36517,  63, 142898)
)

foreach ($row in $rows) {
    $r = $row[0]

    # Keep the date column as plain text (matching the rest of column A)
    # instead of letting it be auto-recognised as a date value.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[1]

    $ws.Cells.Item($r, 2).Value = "overview"
    $ws.Cells.Item($r, 3).Value = "K02000001"
    $ws.Cells.Item($r, 4).Value = "United Kingdom"
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 8).Value = $row[5]
}

Write-Output "Added rows 454-460 to covid_totals"
